# Applies the "Criterion 1 / Criterion 2" sheet swap described by the diff:
#   - Swap the names of the "Criterion 2, Air Speed 0.1" and
#     "Criterion 1, Air Speed 0.1" worksheets (the data stays attached to
#     the worksheet it physically lives on; the *label* moves instead).
#   - Because the label moves rather than the data, the numeric
#     IES/MF results columns (C:D) on those two sheets must be swapped so
#     that the sheet now called "Criterion 1" shows the criterion-1 numbers
#     (and vice-versa).
#   - Two rows (21 and 24) differ by whether a trailing zero "Relative
#     Change" cell (column F) physically exists; that quirk must follow the
#     data, so it is replicated explicitly with ClearContents / assignment.
#   - The "readme" index sheet's table is reshuffled from
#     (index, sheet_name, Date, JobNo, Author) to
#     (index, Author, JobNo, sheet_name, Date), the JobNo/Author values are
#     carried over as-is, and the Date is bumped to the new run date.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Grab stable references to the two worksheets BEFORE renaming so later
#    code can't get confused about which physical sheet is which.
# ---------------------------------------------------------------------
$wsWasCrit2 = $wb.Worksheets.Item("Criterion 2, Air Speed 0.1")   # rId3 / sheet3.xml -> becomes "Criterion 1"
$wsWasCrit1 = $wb.Worksheets.Item("Criterion 1, Air Speed 0.1")   # rId5 / sheet5.xml -> becomes "Criterion 2"

# ---------------------------------------------------------------------
# 2. Swap the numeric results (IES Results / MF Results, columns C:D,
#    rows 2-32) between the two sheets.
# ---------------------------------------------------------------------
$rangeWasCrit2 = $wsWasCrit2.Range("C2:D32")
$rangeWasCrit1 = $wsWasCrit1.Range("C2:D32")

$valuesWasCrit2 = $rangeWasCrit2.Value2
$valuesWasCrit1 = $rangeWasCrit1.Value2

$rangeWasCrit2.Value2 = $valuesWasCrit1
$rangeWasCrit1.Value2 = $valuesWasCrit2

# ---------------------------------------------------------------------
# 3. Row 21 and row 24 also swap whether the trailing (zero) column-F cell
#    is physically present. The sheet that used to be "Criterion 2" loses
#    its F21/F24 cells; the sheet that used to be "Criterion 1" gains them.
# ---------------------------------------------------------------------
$wsWasCrit2.Range("F21").ClearContents()
$wsWasCrit2.Range("F24").ClearContents()

$wsWasCrit1.Range("F21").Value2 = 0
$wsWasCrit1.Range("F24").Value2 = 0

# ---------------------------------------------------------------------
# 4. Swap the sheet names/tabs themselves (via a temporary name so the two
#    assignments don't collide).
# ---------------------------------------------------------------------
$wsWasCrit2.Name = "__tmp_swap__"
$wsWasCrit1.Name = "Criterion 2, Air Speed 0.1"
$wsWasCrit2.Name = "Criterion 1, Air Speed 0.1"

# ---------------------------------------------------------------------
# 5. Update table column headers so they track the sheet's new identity.
#    (Setting the header cell text directly, rather than ListColumn.Name,
#    so the table definition's column name is actually updated.)
# ---------------------------------------------------------------------
$wsWasCrit2.Range("E1").Value2 = "Criterion 1 Absolute Change"
$wsWasCrit2.Range("F1").Value2 = "Criterion 1 Relative Change (%)"

$wsWasCrit1.Range("E1").Value2 = "Criterion 2 Absolute Change"
$wsWasCrit1.Range("F1").Value2 = "Criterion 2 Relative Change (%)"

# ---------------------------------------------------------------------
# 6. Reshuffle the "readme" index sheet: columns go from
#    (index, sheet_name, Date, JobNo, Author) to
#    (index, Author, JobNo, sheet_name, Date), and the run date changes.
#    The sheet_name values are re-derived from the (now renamed) tabs, in
#    the same top-to-bottom tab order as before, so row 3/row 5 correctly
#    pick up the swapped "Criterion 1"/"Criterion 2" labels.
# ---------------------------------------------------------------------
$readme = $wb.Worksheets.Item("readme")

$oldJobNo = $readme.Range("D2").Value2
$oldAuthor = $readme.Range("E2").Value2

$newSheetNames = @(
    $wb.Worksheets.Item("Criteria Failing, Air Speed 0.1").Name,
    $wb.Worksheets.Item(3).Name,
    $wb.Worksheets.Item(4).Name,
    $wb.Worksheets.Item(5).Name
)

# Header row.
$readme.Range("B1").Value2 = "Author"
$readme.Range("C1").Value2 = "JobNo"
$readme.Range("D1").Value2 = "sheet_name"
$readme.Range("E1").Value2 = "Date"

# Data rows: Author, JobNo, sheet_name, Date (sheet_name now reflects the
# swapped tab names; Date is the new run date).
for ($i = 0; $i -lt 4; $i++) {
    $r = 2 + $i
    $readme.Range("B$r").Value2 = $oldAuthor
    $readme.Range("C$r").Value2 = $oldJobNo
    $readme.Range("D$r").Value2 = $newSheetNames[$i]
    $readme.Range("E$r").Value2 = "20220422"
}

Write-Output "done"
